$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.921.06"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.845.46"
$ws.Range("E3").Value = "  +1.17%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "706.11"
$ws.Range("E5").Value = "  +0.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.60"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.843.78"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.65%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.67"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.495.26"
$ws.Range("E15").Value = "  +1.24%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.851.45"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.966.66"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.35"
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "495.84"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.62"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.716"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.49"
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.65"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.19"
$ws.Range("E27").Value = "  -2.10%  "
$ws.Range("E28").Value = "  -3.33%  "
$ws.Range("E29").Value = "  +2.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.41"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  -2.52%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.18"
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.803.37"
$ws.Range("E36").Value = "  +1.38%  "
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.36"
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.03"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("E41").Value = "  +5.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("E42").Value = "  -3.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000318"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.47"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.66"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "415.49"
$ws.Range("E49").Value = "  +1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.298"
$ws.Range("E50").Value = "  -1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.61"
$ws.Range("E51").Value = "  +0.26%  "
